# Apply the "contingencies with rene fine" edit:
#  - extend the header row (row 1) with two new columns P (14) and Q (15),
#    matching the bold/centered/bordered style used by the rest of row 1
#  - for every data row (2-25) flip columns I/K/M/O (1<->2) and append
#    two new data columns P and Q, both valued 2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: new columns P1, Q1 ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
# Copy formatting (bold, centered, bordered) from the existing O1 header cell
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# --- Data rows 2-25 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new column
    $ws.Cells.Item($r, 17).Value = 2   # Q: new column
}
